$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Try exact white (FFFFFF) - should that map to theme 0 without tint?
$ws.Range("A30").Interior.Color = 16777215
$ws.Range("A30").Value = "white"
